$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix C44 data correction: Bakeoff DiFi -> Isolation
$ws.Range("C44").Value = "Isolation"

# Copy style (s=1) from A2 onto the new D column range D2:D89 before filling values
$ws.Range("A2").Copy()
$ws.Range("D2:D89").PasteSpecial(-4122)

# D1 header (no special style, like A1/B1/C1)
$ws.Range("D1").Value = "Subset"

# Fill D2:D89 in an order that reproduces original shared-string insertion order
$ws.Range("D2").Value = "Untreated"
$ws.Range("D3").Value = "DiFi Supermere Day 0"
$ws.Range("D5").Value = "DiFi Supermere Day 4"
$ws.Range("D4").Value = "DiFi Supermere Day 2"
$ws.Range("D6").Value = "Untreated"
$ws.Range("D7").Value = "DiFi Supermere Day 0"
$ws.Range("D8").Value = "DiFi Supermere Day 2"
$ws.Range("D9").Value = "DiFi Supermere Day 4"
$ws.Range("D10").Value = "Untreated"
$ws.Range("D11").Value = "DiFi Supermere Day 0"
$ws.Range("D12").Value = "DiFi Supermere Day 2"
$ws.Range("D13").Value = "DiFi Supermere Day 4"
$ws.Range("D14").Value = "Untreated"
$ws.Range("D15").Value = "CC-CR Supermere Day 0"
$ws.Range("D16").Value = "CC-CR Supermere Day 2"
$ws.Range("D17").Value = "CC-CR Supermere Day 4"
$ws.Range("D18").Value = "Untreated"
$ws.Range("D19").Value = "CC-CR Supermere Day 0"
$ws.Range("D20").Value = "CC-CR Supermere Day 2"
$ws.Range("D21").Value = "CC-CR Supermere Day 4"
$ws.Range("D22").Value = "Untreated"
$ws.Range("D23").Value = "CC-CR Supermere Day 0"
$ws.Range("D24").Value = "CC-CR Supermere Day 2"
$ws.Range("D25").Value = "CC-CR Supermere Day 4"
$ws.Range("D26").Value = "Untreated"
$ws.Range("D27").Value = "DiFi Exomere Day 0"
$ws.Range("D28").Value = "DiFi Exomere Day 2"
$ws.Range("D29").Value = "DiFi Exomere Day 4"
$ws.Range("D30").Value = "Untreated"
$ws.Range("D31").Value = "DiFi Exomere Day 0"
$ws.Range("D32").Value = "DiFi Exomere Day 2"
$ws.Range("D33").Value = "DiFi Exomere Day 4"
$ws.Range("D34").Value = "Untreated"
$ws.Range("D35").Value = "Supermere INPUT"
$ws.Range("D36").Value = "Supermere MOCK"
$ws.Range("D37").Value = "Supermere RNAse temp"
$ws.Range("D38").Value = "RNAse temp"
$ws.Range("D39").Value = "RNAse MOCK"
$ws.Range("D40").Value = "Untreated"
$ws.Range("D41").Value = "DiFi Supermere"
$ws.Range("D42").Value = "DiFi Supermere + TGFBI"
$ws.Range("D43").Value = "DiFi Supermere TGFBI OE"
$ws.Range("D44").Value = "TGFBI"
$ws.Range("D45").Value = "Untreated"
$ws.Range("D46").Value = "DiFi FPLC Supermere"
$ws.Range("D47").Value = "DiFi FPLC Exomere"
$ws.Range("D48").Value = "DiFi UC Supermere"
$ws.Range("D49").Value = "DiFi UC Exomere"
$ws.Range("D50").Value = "DiFi UC sEV"
$ws.Range("D51").Value = "Untreated"
$ws.Range("D52").Value = "CC-CR FPLC Supermere"
$ws.Range("D53").Value = "CC-CR FPLC Exomere"
$ws.Range("D54").Value = "CC-CR UC Supermere"
$ws.Range("D55").Value = "CC-CR UC Exomere"
$ws.Range("D56").Value = "CC-CR UC sEV"
$ws.Range("D57").Value = "Untreated"
$ws.Range("D58").Value = "Untreated"
$ws.Range("D59").Value = "Untreated"
$ws.Range("D60").Value = "LPS"
$ws.Range("D61").Value = "LPS"
$ws.Range("D62").Value = "LPS"
$ws.Range("D63").Value = "Healthy"
$ws.Range("D64").Value = "Healthy"
$ws.Range("D65").Value = "Healthy"
$ws.Range("D66").Value = "Healthy"
$ws.Range("D67").Value = "Healthy"
$ws.Range("D68").Value = "Healthy"
$ws.Range("D69").Value = "Healthy"
$ws.Range("D70").Value = "Healthy"
$ws.Range("D71").Value = "Healthy"
$ws.Range("D72").Value = "CRC"
$ws.Range("D73").Value = "CRC"
$ws.Range("D74").Value = "CRC"
$ws.Range("D75").Value = "CRC"
$ws.Range("D76").Value = "CRC"
$ws.Range("D77").Value = "CRC"
$ws.Range("D78").Value = "CRC"
$ws.Range("D79").Value = "CRC"
$ws.Range("D80").Value = "CRC"
$ws.Range("D81").Value = "Healthy no cells"
$ws.Range("D82").Value = "Healthy no cells"
$ws.Range("D83").Value = "Healthy no cells"
$ws.Range("D84").Value = "CRC no cells"
$ws.Range("D85").Value = "CRC no cells"
$ws.Range("D86").Value = "CRC no cells"
$ws.Range("D87").Value = "Plasma"
$ws.Range("D88").Value = "Plasma"
$ws.Range("D89").Value = "Plasma"
